$wb = $excel.ActiveWorkbook

# Sheet "sum_response_time": reduce REST (column B) values for rest2 and rest4
$ws1 = $wb.Worksheets.Item("sum_response_time")
$ws1.Range("B3").Value = 18.76
$ws1.Range("B5").Value = 16.9

# Sheet "total_data_transferred": reduce REST (column B) values for rest2 and rest4
$ws2 = $wb.Worksheets.Item("total_data_transferred")
$ws2.Range("B3").Value = 9.3271484375
$ws2.Range("B5").Value = 1.1728515625
